$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "N splits" column header (C1) to "Start", and add a new "End" header (D1)
$ws.Range("C1").Value = "Start"
$ws.Range("D1").Value = "End"

# Row 2: rename the templated output path to the first explicit split output,
# and add its End time (column D). Input (A2) and N-splits count (C2) stay the same.
$ws.Range("B2").Value = "out/test2/part1.mp4"
$ws.Range("D2").Value = 12

# Row 3: new second split entry (Input / Output / Start / End)
$ws.Range("A3").Value = "in/test.mp4"
$ws.Range("B3").Value = "out/test2/part2.mp4"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 15

# Restore the selection to match the saved view state
$ws.Range("B5").Select()
